$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Datos actualizados" timestamp refresh (14:30 -> 15:47)
$ws.Range("A1").Value = "Datos actualizados a 20 de Agosto de 2020 a las 15:47"

# Paises Bajos overtakes Guatemala in the ranking (rows 45/46 swap countries)
$ws.Range("A45").Value = "Paises Bajos"
$ws.Range("A46").Value = "Guatemala"

# Montserrat overtakes Islas Malvinas in the ranking (rows 213/214 swap countries)
$ws.Range("A213").Value = "Montserrat"
$ws.Range("A214").Value = "Islas Malvinas"

# Refresh Casos totales / Nuevos casos / Casos activos / Recuperados / Casos criticos / Muertes hoy / Muertes
# Row 4
$ws.Range("B4").Value = 5702782
$ws.Range("C4").Value = 1851
$ws.Range("D4").Value = 3063259
$ws.Range("E4").Value = 2463128
$ws.Range("G4").Value = 58
$ws.Range("H4").Value = 176395

# Row 6
$ws.Range("B6").Value = 2858346
$ws.Range("C6").Value = 22524
$ws.Range("D6").Value = 2115054
$ws.Range("E6").Value = 689047
$ws.Range("G6").Value = 251
$ws.Range("H6").Value = 54245

# Row 17
$ws.Range("B17").Value = 303973
$ws.Range("C17").Value = 1287
$ws.Range("D17").Value = 275476
$ws.Range("E17").Value = 24949
$ws.Range("G17").Value = 42
$ws.Range("H17").Value = 3548

# Row 24
$ws.Range("B24").Value = 192797
$ws.Range("C24").Value = 3995
$ws.Range("D24").Value = 137200
$ws.Range("E24").Value = 49389
$ws.Range("G24").Value = 87
$ws.Range("H24").Value = 6208

# Row 28
$ws.Range("B28").Value = 116224
$ws.Range("C28").Value = 268
$ws.Range("D28").Value = 112924
$ws.Range("E28").Value = 3107

# Row 45
$ws.Range("B45").Value = 65054
$ws.Range("C45").Value = 529
$ws.Range("D45").Value = 0
$ws.Range("E45").Value = 0
$ws.Range("G45").Value = 10
$ws.Range("H45").Value = 6191

# Row 46
$ws.Range("B46").Value = 64881
$ws.Range("D46").Value = 53362
$ws.Range("E46").Value = 9052
$ws.Range("H46").Value = 2467

# Row 50
$ws.Range("B50").Value = 54992
$ws.Range("C50").Value = 291
$ws.Range("D50").Value = 40264
$ws.Range("E50").Value = 12940
$ws.Range("G50").Value = 2
$ws.Range("H50").Value = 1788

# Row 55
$ws.Range("B55").Value = 43260
$ws.Range("C55").Value = 166
$ws.Range("D55").Value = 41276
$ws.Range("E55").Value = 1723
$ws.Range("G55").Value = 5
$ws.Range("H55").Value = 261

# Row 63
$ws.Range("B63").Value = 34759
$ws.Range("C63").Value = 139
$ws.Range("D63").Value = 32511
$ws.Range("E63").Value = 1738

# Row 68
$ws.Range("B68").Value = 30209
$ws.Range("C68").Value = 161
$ws.Range("D68").Value = 27908
$ws.Range("E68").Value = 1612
$ws.Range("G68").Value = 5
$ws.Range("H68").Value = 689

# Row 162
$ws.Range("B162").Value = 887
$ws.Range("C162").Value = 2
$ws.Range("D162").Value = 829
$ws.Range("E162").Value = 43

# Row 173
$ws.Range("B173").Value = 399
$ws.Range("C173").Value = 5
$ws.Range("E173").Value = 60

# Row 213
$ws.Range("D213").Value = 12
$ws.Range("H213").Value = 1

# Row 214
$ws.Range("D214").Value = 13
$ws.Range("H214").Value = 0
